# Add four new product rows (122-125) to the price list, mirroring the
# scraped-data rows already in the sheet. Per the commit message ("alterando
# no format data para float o preco"), the two most recently captured prices
# (rows 124-125) are stored as real numeric/float values, while the other two
# retain their raw, messily-formatted text (exactly as scraped) like most of
# the pre-existing rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 122 - price kept as literal text "79.99"
$ws.Range("A122").Value = 1567227175
$ws.Range("B122").Value = "Abraçadeiras De Nylon Para Lacre Brancas 7 6mm X 500mm"
# Force text storage so "79.99" isn't auto-coerced into a number, then drop
# the number-format override so the cell keeps the sheet's default style.
$ws.Range("C122").NumberFormat = "@"
$ws.Range("C122").Value = "79.99"
$ws.Range("C122").ClearFormats()

# Row 123 - price kept as literal text "2.399.00"
$ws.Range("A123").Value = 89801243
$ws.Range("B123").Value = "Máquina de Pintura Elétrica 900W Airless MPA 120 220V Vonder"
$ws.Range("C123").Value = "2.399.00"

# Row 124 - price stored as a real float
$ws.Range("A124").Value = 90795621
$ws.Range("B124").Value = "Furadeira e Parafusadeira de Impacto a Bateria Bosch com Carregador e Bateria 18V 1/2`" GSB 180-Li Bivolt"
$ws.Range("C124").Value = 619.9

# Row 125 - price stored as a real float
$ws.Range("A125").Value = 89837055
$ws.Range("B125").Value = "Carvão Briquete 2,5Kg Pérola Negra"
$ws.Range("C125").Value = 18.29
